# February 11 - commit
#
# On the "ManageOrdersPage" sheet (the active/selected sheet), the demo
# "Time Field 1 / Time Field 2 / Expected Text" sample columns (B:D) are
# removed from the header row and the first data row - the cells are
# blanked out but keep their existing styles. Row 2 shrinks back down to
# the normal row height now that the tall wrapped text is gone, and the
# selection moves to A6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ManageOrdersPage")

# Blank out the header cells B1:D1 (style stays s="3").
$ws.Range("B1:D1").ClearContents()

# Blank out the matching data cells B2:D2 (styles stay s="1"/"7"/"8").
$ws.Range("B2:D2").ClearContents()

# Row 2 no longer needs the taller wrapped-text height - match the other
# data rows (18pt).
$ws.Rows.Item(2).RowHeight = 18

# Move the active selection to A6.
$ws.Range("A6").Select() | Out-Null
